$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. Remove existing hyperlinks up front. This COM shim does not shift a
#    Hyperlink's anchor Range when rows/columns are inserted or deleted, so
#    we manage hyperlinks manually and re-add them once the grid is final.
# ---------------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 1. Drop the "E80H30" product row (old row 4) and pull the "E40H8..." row
#    (old row 5) up into its place, leaving old row 5 blank (matching the
#    row-5/row-4 style pattern, which were identical before this edit).
#    We copy only A:F (not G) because G holds the helper "0" value that must
#    stay put in both rows.
# ---------------------------------------------------------------------------
$ws.Range("A5:F5").Copy()
$ws.Range("A4").PasteSpecial(-4104)   # xlPasteAll
$ws.Range("A5:F5").ClearContents()
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Insert a new column at E ("Çıkış Tipi" / Output type). Existing E:H
#    shift right to F:I.
# ---------------------------------------------------------------------------
$ws.Columns("E:E").Insert()

# ---------------------------------------------------------------------------
# 3. Populate the new column E. Copy D's formatting down so styles match
#    (header style for E1, data style for E2:E19).
# ---------------------------------------------------------------------------
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("E1").Value = "Çıkış Tipi"

$ws.Range("D2:D19").Copy()
$ws.Range("E2:E19").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("E2").Value = "Push Pull"
$ws.Range("E3").Value = "TTL"
$ws.Range("E4").Value = "Line Driver"

# ---------------------------------------------------------------------------
# 4. New note cell next to the ARC H 50 row.
# ---------------------------------------------------------------------------
$ws.Range("I3").Value = "ARC-H-50-3600-HTL-6-3M-FZ"

# ---------------------------------------------------------------------------
# 5. Re-create the hyperlinks at their final locations (G2:G4).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("G2"), "http://www.fenac.com.tr/tr/tam-delik-artimli-rotary-enkoderler/20-o100-fnc-100h-serisi-tam-delikli-artimli-rotary-enkoder.html", "/besleme_gerilimi-5_30vdc/k_devresi-abztersleri/balant_ekli-2mt_kablolu/gvde_ap-100mm/delik_ap-23mm/puls-5000 ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "http://www.alfasanayi.com/arc-h-50-hollow-saft-atek-made-in-turkey-urun-176.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G4"), "http://www.endelkon.com/autonics-e40h8-25006-l-5-pmu1553") | Out-Null

# ---------------------------------------------------------------------------
# 6. Column widths (character units). This runtime adds a fixed 5/6 padding
#    and rounds to the nearest 1/6 when a ColumnWidth is assigned, so we
#    back-solve for the input that lands closest to the desired stored width.
#    Column G keeps its original (pre-edit) width explicitly, since the
#    insert operation otherwise drags the old column F's width into G.
# ---------------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 15.5
$ws.Columns("B:B").ColumnWidth = 7.16666666666667
$ws.Columns("D:D").ColumnWidth = 10.16666666666667
$ws.Columns("E:E").ColumnWidth = 13.6666666666667
$ws.Columns("F:F").ColumnWidth = 4.83333333333333
$ws.Columns("G:G").ColumnWidth = 16.6666666666667
$ws.Columns("I:I").ColumnWidth = 25.1666666666667

# ---------------------------------------------------------------------------
# 7. Selection.
# ---------------------------------------------------------------------------
$ws.Range("G3").Select()
